$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Accredited)
$ws.Range("D2").Value = 268
$ws.Range("E2").Value = 106
$ws.Range("F2").Value = 105
$ws.Range("G2").Value = 159
$ws.Range("H2").Value = 131
$ws.Range("K2").Value = 236
$ws.Range("L2").Value = 204
$ws.Range("N2").Value = 140

# Row 3 (Unaccredited)
$ws.Range("C3").Value = 66
$ws.Range("D3").Value = 354
$ws.Range("E3").Value = 153
$ws.Range("F3").Value = 171
$ws.Range("G3").Value = 255
$ws.Range("H3").Value = 186
$ws.Range("I3").Value = 73
$ws.Range("J3").Value = 186
$ws.Range("K3").Value = 353
$ws.Range("L3").Value = 358
$ws.Range("M3").Value = 157
$ws.Range("N3").Value = 147
